$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value "computer" in cell D3 (new column D, row 3)
$ws.Range("D3").Value = "computer"

# Update the selected cell to match the new active cell D3
$ws.Range("D3").Select()
